# Update market-price derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/WVR sheets, mirroring a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 46127
$ws.Cells.Item(87, 10).Value = 46127
$ws.Cells.Item(87, 12).Value = 46127
$ws.Cells.Item(87, 14).Value = -48623

$ws.Cells.Item(90, 8).Value = 46127
$ws.Cells.Item(90, 10).Value = 46127
$ws.Cells.Item(90, 12).Value = 138381
$ws.Cells.Item(90, 14).Value = -150861

$ws.Cells.Item(97, 8).Value = 480
$ws.Cells.Item(97, 9).Value = 480
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1440
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -944
$ws.Cells.Item(97, 14).ClearContents() | Out-Null

$ws.Cells.Item(103, 8).Value = 653.4211
$ws.Cells.Item(103, 9).Value = 642
$ws.Cells.Item(103, 10).Value = 696.25
$ws.Cells.Item(103, 11).Value = 1926
$ws.Cells.Item(103, 12).Value = 2088.75
$ws.Cells.Item(103, 13).Value = -1340
$ws.Cells.Item(103, 14).Value = -3260.75

$ws.Cells.Item(112, 8).Value = 1405.7084
$ws.Cells.Item(112, 10).Value = 1499.6279
$ws.Cells.Item(112, 12).Value = 4498.8837
$ws.Cells.Item(112, 14).Value = -6714.8837

$ws.Cells.Item(115, 8).Value = 8333973
$ws.Cells.Item(115, 9).Value = 8333973
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 25001919
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = -25000352
$ws.Cells.Item(115, 14).ClearContents() | Out-Null

$ws.Cells.Item(118, 8).Value = 422.6
$ws.Cells.Item(118, 9).Value = 298.33334
$ws.Cells.Item(118, 10).Value = 609
$ws.Cells.Item(118, 11).Value = 895.0000200000001
$ws.Cells.Item(118, 12).Value = 1827
$ws.Cells.Item(118, 13).Value = 761.9999799999999
$ws.Cells.Item(118, 14).Value = -5141

$ws.Cells.Item(123, 8).Value = 66390
$ws.Cells.Item(123, 10).Value = 66390
$ws.Cells.Item(123, 12).Value = 66390
$ws.Cells.Item(123, 14).Value = -76190

$ws.Cells.Item(124, 8).Value = 60319.332
$ws.Cells.Item(124, 10).Value = 60319.332
$ws.Cells.Item(124, 12).Value = 60319.332
$ws.Cells.Item(124, 14).Value = -70139.33199999999

$ws.Cells.Item(126, 8).Value = 69660
$ws.Cells.Item(126, 10).Value = 69660
$ws.Cells.Item(126, 12).Value = 69660
$ws.Cells.Item(126, 14).Value = -79540

$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 14).ClearContents() | Out-Null

$ws.Cells.Item(130, 8).Value = 20926.666
$ws.Cells.Item(130, 10).Value = 42780
$ws.Cells.Item(130, 12).Value = 42780
$ws.Cells.Item(130, 14).Value = -52820

$ws.Cells.Item(137, 8).Value = 239152.55
$ws.Cells.Item(137, 9).Value = 362086.28
$ws.Cells.Item(137, 10).Value = 2741.5386
$ws.Cells.Item(137, 11).Value = 1086258.84
$ws.Cells.Item(137, 12).Value = 8224.6158
$ws.Cells.Item(137, 13).Value = -1083708.84
$ws.Cells.Item(137, 14).Value = -13324.6158

$ws.Cells.Item(138, 8).Value = 3931.8357
$ws.Cells.Item(138, 9).Value = 1907
$ws.Cells.Item(138, 10).Value = 4370.55
$ws.Cells.Item(138, 11).Value = 5721
$ws.Cells.Item(138, 12).Value = 13111.65
$ws.Cells.Item(138, 13).Value = -581
$ws.Cells.Item(138, 14).Value = -23391.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 90
$ws.Cells.Item(5, 9).Value = 90
$ws.Cells.Item(5, 11).Value = 90
$ws.Cells.Item(5, 13).Value = 22

$ws.Cells.Item(74, 8).Value = 1247.84
$ws.Cells.Item(74, 9).Value = 718.8421
$ws.Cells.Item(74, 10).Value = 2923
$ws.Cells.Item(74, 11).Value = 718.8421
$ws.Cells.Item(74, 12).Value = 2923
$ws.Cells.Item(74, 13).Value = 155.1579
$ws.Cells.Item(74, 14).Value = -4671

$ws.Cells.Item(77, 8).Value = 1247.84
$ws.Cells.Item(77, 9).Value = 718.8421
$ws.Cells.Item(77, 10).Value = 2923
$ws.Cells.Item(77, 11).Value = 3594.2105
$ws.Cells.Item(77, 12).Value = 14615
$ws.Cells.Item(77, 13).Value = 773.7895000000003
$ws.Cells.Item(77, 14).Value = -23351

$ws.Cells.Item(102, 8).Value = 3165
$ws.Cells.Item(102, 9).Value = 2221.6667
$ws.Cells.Item(102, 11).Value = 2221.6667
$ws.Cells.Item(102, 13).Value = -599.6667000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 90
$ws.Cells.Item(4, 9).Value = 90
$ws.Cells.Item(4, 11).Value = 90
$ws.Cells.Item(4, 13).Value = 25

$ws.Cells.Item(105, 8).Value = 2112.75
$ws.Cells.Item(105, 9).Value = 1353.3334
$ws.Cells.Item(105, 10).Value = 2872.1667
$ws.Cells.Item(105, 11).Value = 1353.3334
$ws.Cells.Item(105, 12).Value = 2872.1667
$ws.Cells.Item(105, 13).Value = 393.6666
$ws.Cells.Item(105, 14).Value = -6366.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 47
$ws.Cells.Item(7, 9).Value = 59.5
$ws.Cells.Item(7, 10).Value = 22
$ws.Cells.Item(7, 11).Value = 59.5
$ws.Cells.Item(7, 12).Value = 22
$ws.Cells.Item(7, 13).Value = 53.5
$ws.Cells.Item(7, 14).Value = -248

$ws.Cells.Item(17, 8).Value = 22500
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 22500
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 22500
$ws.Cells.Item(17, 13).ClearContents() | Out-Null
$ws.Cells.Item(17, 14).Value = -22848

$ws.Cells.Item(25, 8).Value = 11200
$ws.Cells.Item(25, 10).Value = 11200
$ws.Cells.Item(25, 12).Value = 11200
$ws.Cells.Item(25, 14).Value = -11548

$ws.Cells.Item(58, 8).Value = 1067.9756
$ws.Cells.Item(58, 9).Value = 793
$ws.Cells.Item(58, 11).Value = 793
$ws.Cells.Item(58, 13).Value = -590

$ws.Cells.Item(134, 8).Value = 2891.913
$ws.Cells.Item(134, 9).Value = 1276.1111
$ws.Cells.Item(134, 10).Value = 3930.6428
$ws.Cells.Item(134, 11).Value = 3828.3333
$ws.Cells.Item(134, 12).Value = 11791.9284
$ws.Cells.Item(134, 13).Value = -1293.3333
$ws.Cells.Item(134, 14).Value = -16861.9284

$ws.Cells.Item(136, 8).Value = 1067.9756
$ws.Cells.Item(136, 9).Value = 793
$ws.Cells.Item(136, 11).Value = 2379
$ws.Cells.Item(136, 13).Value = 171

$ws.Cells.Item(141, 8).Value = 56138.555
$ws.Cells.Item(141, 10).Value = 56138.555
$ws.Cells.Item(141, 12).Value = 56138.555
$ws.Cells.Item(141, 14).Value = -66498.55499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 3929356.5
$ws.Cells.Item(132, 9).Value = 872.5333000000001
$ws.Cells.Item(132, 10).Value = 5893598.5
$ws.Cells.Item(132, 11).Value = 7852.7997
$ws.Cells.Item(132, 12).Value = 53042386.5
$ws.Cells.Item(132, 13).Value = -5322.7997
$ws.Cells.Item(132, 14).Value = -53047446.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 14999.75
$ws.Cells.Item(54, 10).Value = 14999.75
$ws.Cells.Item(54, 12).Value = 14999.75
$ws.Cells.Item(54, 14).Value = -16039.75

$ws.Cells.Item(136, 8).Value = 4877.815
$ws.Cells.Item(136, 9).Value = 5329.136
$ws.Cells.Item(136, 10).Value = 2892
$ws.Cells.Item(136, 11).Value = 15987.408
$ws.Cells.Item(136, 12).Value = 8676
$ws.Cells.Item(136, 13).Value = -13437.408
$ws.Cells.Item(136, 14).Value = -13776
